$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of row 27 down to the new row 28 before filling values,
# so the new row picks up the same cell styles (centered / wrap / date format)
# without introducing new style records.
$ws.Range("A27:I27").Copy()
$ws.Range("A28:I28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the tag string on row 27 (problem 875, Koko Eating Bananas):
# "#two-pointers #array #binary-search #必背 " -> "#two-pointers #array #binary-search #核心 "
$ws.Range("C27").Value = "#two-pointers #array #binary-search #核心 "

# Add new row 28 for problem 1283: Find the Smallest Divisor Given a Threshold
$ws.Range("A28").Value = 1283
$ws.Range("B28").Value = "Find the Smallest Divisor Given a Threshold"
$ws.Range("C28").Value = "#two-pointers #array #binary-search #核心 "
$ws.Range("D28").Value = "medium"
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = 45838
$ws.Range("I28").Value = 45838

$ws.Rows.Item(28).RowHeight = 51

$ws.Range("I28").Select()
$excel.ActiveWindow.ScrollRow = 24
